$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1103.3846
$ws.Range("I53").Value = 84.125
$ws.Range("J53").Value = 2734.2
$ws.Range("K53").Value = 84.125
$ws.Range("L53").Value = 2734.2
$ws.Range("M53").Value = 552.875
$ws.Range("N53").Value = -4008.2
$ws.Range("H64").Value = 5438.357
$ws.Range("I64").Value = 5449
$ws.Range("K64").Value = 5449
$ws.Range("M64").Value = -5201
$ws.Range("H67").Value = 5438.357
$ws.Range("I67").Value = 5449
$ws.Range("K67").Value = 5449
$ws.Range("M67").Value = -4591
$ws.Range("H113").Value = 61030.06
$ws.Range("I113").Value = 2455.4443
$ws.Range("J113").Value = 126926.5
$ws.Range("K113").Value = 2455.4443
$ws.Range("L113").Value = 126926.5
$ws.Range("M113").Value = 798.5556999999999
$ws.Range("N113").Value = -133434.5
$ws.Range("H116").Value = 7925.6665
$ws.Range("J116").Value = 8666.333000000001
$ws.Range("L116").Value = 8666.333000000001
$ws.Range("N116").Value = -15550.333
$ws.Range("H132").Value = 1572.0714
$ws.Range("I132").Value = 1187.1818
$ws.Range("J132").Value = 2983.3333
$ws.Range("K132").Value = 3561.5454
$ws.Range("L132").Value = 8949.999899999999
$ws.Range("M132").Value = -1031.5454
$ws.Range("N132").Value = -14009.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4893.2
$ws.Range("I32").Value = 5155.515
$ws.Range("K32").Value = 5155.515
$ws.Range("M32").Value = -4868.515
$ws.Range("H36").Value = 6499.6665
$ws.Range("I36").Value = 6499.6665
$ws.Range("K36").Value = 6499.6665
$ws.Range("M36").Value = -6153.6665
$ws.Range("H45").Value = 4127.5
$ws.Range("I45").Value = 3931.4285
$ws.Range("K45").Value = 3931.4285
$ws.Range("M45").Value = -3554.4285
$ws.Range("H61").Value = 4858.75
$ws.Range("I61").Value = 4002.077
$ws.Range("K61").Value = 4002.077
$ws.Range("M61").Value = -3790.077
$ws.Range("H102").Value = 3264
$ws.Range("I102").Value = 3264
$ws.Range("K102").Value = 3264
$ws.Range("M102").Value = -1642
$ws.Range("H110").Value = 18751622
$ws.Range("I110").Value = 21429944
$ws.Range("J110").Value = 3365.6667
$ws.Range("K110").Value = 21429944
$ws.Range("L110").Value = 3365.6667
$ws.Range("M110").Value = -21427899
$ws.Range("N110").Value = -7455.6667
$ws.Range("H132").Value = 6632.7666
$ws.Range("I132").Value = 3955.8696
$ws.Range("K132").Value = 11867.6088
$ws.Range("M132").Value = -9337.6088
$ws.Range("H136").Value = 4858.75
$ws.Range("I136").Value = 4002.077
$ws.Range("K136").Value = 12006.231
$ws.Range("M136").Value = -9456.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 16491.5
$ws.Range("I54").Value = 16491.5
$ws.Range("K54").Value = 16491.5
$ws.Range("M54").Value = -16007.5
$ws.Range("H97").Value = 9690.5
$ws.Range("I97").Value = 9690.5
$ws.Range("K97").Value = 9690.5
$ws.Range("M97").Value = -8699.5
$ws.Range("H105").Value = 55570068
$ws.Range("I105").Value = 66683984
$ws.Range("K105").Value = 66683984
$ws.Range("M105").Value = -66682237

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 390.42856
$ws.Range("I22").Value = 350.6
$ws.Range("K22").Value = 350.6
$ws.Range("M22").Value = -0.6000000000000227
$ws.Range("H31").Value = 62506350
$ws.Range("I31").Value = 1000000000
$ws.Range("K31").Value = 1000000000
$ws.Range("M31").Value = -999999705
$ws.Range("H34").Value = 62506350
$ws.Range("I34").Value = 1000000000
$ws.Range("K34").Value = 1000000000
$ws.Range("M34").Value = -999999798
$ws.Range("H134").Value = 7924.448
$ws.Range("I134").Value = 7864.4
$ws.Range("K134").Value = 23593.2
$ws.Range("M134").Value = -21058.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 446
$ws.Range("I39").Value = 419
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 1257
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -963
$ws.Range("N39").Value = -2088
$ws.Range("H63").Value = 4950
$ws.Range("J63").Value = 4950
$ws.Range("L63").Value = 14850
$ws.Range("N63").Value = -16348
$ws.Range("H66").Value = 4950
$ws.Range("J66").Value = 4950
$ws.Range("L66").Value = 44550
$ws.Range("N66").Value = -52038
$ws.Range("H114").Value = 3266.4
$ws.Range("I114").Value = 325.25
$ws.Range("J114").Value = 15031
$ws.Range("K114").Value = 975.75
$ws.Range("L114").Value = 45093
$ws.Range("M114").Value = 2278.25
$ws.Range("N114").Value = -51601
$ws.Range("H117").Value = 167880.4
$ws.Range("J117").Value = 186506
$ws.Range("L117").Value = 559518
$ws.Range("N117").Value = -566402
$ws.Range("H121").Value = 23811366
$ws.Range("I121").Value = 363.5
$ws.Range("J121").Value = 33335768
$ws.Range("K121").Value = 1090.5
$ws.Range("L121").Value = 100007304
$ws.Range("M121").Value = 219.5
$ws.Range("N121").Value = -100009924
$ws.Range("H129").Value = 50001304
$ws.Range("I129").Value = 925.7143
$ws.Range("J129").Value = 166668850
$ws.Range("K129").Value = 2777.1429
$ws.Range("L129").Value = 500006550
$ws.Range("M129").Value = 2222.8571
$ws.Range("N129").Value = -500016550
$ws.Range("H131").Value = 13893136
$ws.Range("I131").Value = 27778654
$ws.Range("K131").Value = 83335962
$ws.Range("M131").Value = -83330922
$ws.Range("H132").Value = 1884
$ws.Range("J132").Value = 2798.6
$ws.Range("L132").Value = 25187.4
$ws.Range("N132").Value = -30247.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3143.0833
$ws.Range("I80").Value = 1113
$ws.Range("J80").Value = 4158.125
$ws.Range("K80").Value = 1113
$ws.Range("L80").Value = 4158.125
$ws.Range("M80").Value = -115
$ws.Range("N80").Value = -6154.125
$ws.Range("H83").Value = 3143.0833
$ws.Range("I83").Value = 1113
$ws.Range("J83").Value = 4158.125
$ws.Range("K83").Value = 5565
$ws.Range("L83").Value = 20790.625
$ws.Range("M83").Value = -573
$ws.Range("N83").Value = -30774.625
$ws.Range("H122").Value = 8066.1763
$ws.Range("J122").Value = 5999
$ws.Range("L122").Value = 17997
$ws.Range("N122").Value = -22897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2208.2778
$ws.Range("I82").Value = 1913.6364
$ws.Range("J82").Value = 2671.2856
$ws.Range("K82").Value = 1913.6364
$ws.Range("L82").Value = 2671.2856
$ws.Range("M82").Value = -1552.6364
$ws.Range("N82").Value = -3393.2856
$ws.Range("H85").Value = 2208.2778
$ws.Range("I85").Value = 1913.6364
$ws.Range("J85").Value = 2671.2856
$ws.Range("K85").Value = 1913.6364
$ws.Range("L85").Value = 2671.2856
$ws.Range("M85").Value = -665.6364000000001
$ws.Range("N85").Value = -5167.2856
$ws.Range("H122").Value = 3504.1052
$ws.Range("I122").Value = 3445.8235
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 10337.4705
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -7887.470499999999
$ws.Range("N122").Value = -16898.5
$ws.Range("H132").Value = 6939.773
$ws.Range("I132").Value = 5986.8237
$ws.Range("J132").Value = 10179.8
$ws.Range("K132").Value = 17960.4711
$ws.Range("L132").Value = 30539.4
$ws.Range("M132").Value = -15430.4711
$ws.Range("N132").Value = -35599.39999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 11926.692
$ws.Range("I52").Value = 6945.375
$ws.Range("J52").Value = 19896.8
$ws.Range("K52").Value = 6945.375
$ws.Range("L52").Value = 19896.8
$ws.Range("M52").Value = -6719.375
$ws.Range("N52").Value = -20348.8
$ws.Range("H61").Value = 12367
$ws.Range("I61").Value = 10492.375
$ws.Range("J61").Value = 17366
$ws.Range("K61").Value = 10492.375
$ws.Range("L61").Value = 17366
$ws.Range("M61").Value = -10200.375
$ws.Range("N61").Value = -17950
$ws.Range("H81").Value = 7574
$ws.Range("I81").Value = 3118.4
$ws.Range("K81").Value = 6236.8
$ws.Range("M81").Value = -5175.8
$ws.Range("H84").Value = 7574
$ws.Range("I84").Value = 3118.4
$ws.Range("K84").Value = 31184
$ws.Range("M84").Value = -25880
$ws.Range("H122").Value = 3267.2942
$ws.Range("I122").Value = 3103
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 9309
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -6859
$ws.Range("N122").Value = -18398.5
